# Applies the "Updated files and tables" edit to the worksheet:
#  - adds a new "colors" header label in I8 (bold, centered)
#  - adds/updates boolean flags in column I for rows 9-23
#  - updates the view selection to the cell last edited (L14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "colors" column header in I8, bold + centered ---
$i8 = $ws.Range("I8")
$i8.Value = "colors"
$i8.Style = "Normal"
$i8.Font.Bold = $true
$i8.HorizontalAlignment = -4108  # xlCenter

# --- Reset any inherited formatting on I9:I23 so the new boolean cells ---
# --- pick up the sheet's default (unstyled) cell format               ---
foreach ($r in 9..23) {
    $cell = $ws.Cells.Item($r, 9)  # column I = 9
    $cell.Style = "Normal"
}

# --- New TRUE/FALSE "colors" flag values for each gene row ---
$ws.Range("I9").Value  = $false
$ws.Range("I10").Value = $false
$ws.Range("I11").Value = $true
$ws.Range("I12").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("I14").Value = $false
$ws.Range("I15").Value = $false
$ws.Range("I16").Value = $false
$ws.Range("I17").Value = $false
$ws.Range("I18").Value = $false
$ws.Range("I19").Value = $false
$ws.Range("I20").Value = $false
$ws.Range("I21").Value = $false
$ws.Range("I22").Value = $false
$ws.Range("I23").Value = $false

# --- Move the view/selection to reflect where the author last worked ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("L14").Select() | Out-Null
